$wb = $excel.ActiveWorkbook

# --- Sheet "Orders" (sheet1): append new rows 12-21 ---
$ws1 = $wb.Worksheets.Item("Orders")

function Set-TextValue($cell, $value) {
    # Force text storage (matches existing sheet convention where numeric-looking
    # strings like "5", "8", "2" are kept as text, not numbers).
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

$newRows = @(
    @{ Row = 12; A = $null; C = "579_腊梅红_wax red_undefined_1bunch"; F = "5" },
    @{ Row = 13; A = $null; C = "614_康乃馨绿_green_undefined_20stems"; F = "8" },
    @{ Row = 14; A = $null; C = "611_康乃馨奶油白_cream white_undefined_20stems"; F = "5" },
    @{ Row = 15; A = $null; C = "602_康乃馨白_white_undefined_20stems"; F = "10" },
    @{ Row = 16; A = $null; C = "597_尤加利叶小叶_undefined_undefined_1bunch"; F = "3" },
    @{ Row = 17; A = "2";   C = "138_卡罗拉_Carola_Rosa rugosa Thunb._20stems"; F = "12" },
    @{ Row = 18; A = $null; C = "148_坦尼克_Tineke_Rosa rugosa Thunb._20stems"; F = "15" },
    @{ Row = 19; A = $null; C = "268_猩红泡泡_spray red_Rosa rugosa Thunb._10stems"; F = "6" },
    @{ Row = 20; A = $null; C = "43_拉丝红_Spider Red_Gerbera L._20stems"; F = "5" },
    @{ Row = 21; A = $null; C = "44_拉丝粉_Spider Pink_Gerbera L._20stems"; F = $null }
)

foreach ($r in $newRows) {
    if ($r.A -ne $null) {
        Set-TextValue $ws1.Cells.Item($r.Row, 1) $r.A
    }
    $ws1.Cells.Item($r.Row, 3).Value = $r.C
    if ($r.F -ne $null) {
        Set-TextValue $ws1.Cells.Item($r.Row, 6) $r.F
    }
}

# --- Sheet "Summary" (sheet2): update G2 value ---
$ws2 = $wb.Worksheets.Item("Summary")
Set-TextValue $ws2.Cells.Item(2, 7) "0588103102020555851031215650"
